# Apply the edits described by the diff:
# 1. Rename sheet "Dashboard" -> "DataSet"
# 2. Make "DataSet" (first sheet) the active/selected tab, with cell B25 selected
# 3. Make "FooterLinks" (second sheet) no longer the active/selected tab (its own
#    selection E23 is left untouched)

$wb = $excel.ActiveWorkbook

# Grab the first sheet reference before renaming it
$dashboard = $wb.Worksheets.Item(1)

# 1. Rename "Dashboard" to "DataSet"
$dashboard.Name = "DataSet"

# 2 & 3. Activating the DataSet sheet and selecting B25 makes it the tabSelected /
#        active tab sheet, and consequently FooterLinks loses tabSelected.
$dashboard.Activate()
$dashboard.Range("B25").Select()
